# "added independent variables beyond sex"
#
# On the "sexes" sheet, add two new columns (GrpSize, Age) describing each
# cage beyond the existing CageID/Sex columns, then rename the sheet to
# "Sheet3" and make it the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item("sexes")

# --- New header cells: GrpSize (C1), Age (D1) ---
$ws3.Cells.Item(1, 3).Value = "GrpSize"
$ws3.Cells.Item(1, 4).Value = "Age"
$ws3.Range("A1:D1").Font.Bold = $true

# Row -> GrpSize, Age
$grpSize = @{2=2; 3=2; 4=3; 5=2; 6=2; 7=2; 8=2; 9=2; 10=3; 11=2; 12=3; 13=2}
$age     = @{2="Young"; 3="Young"; 4="Young"; 5="Young"; 6="Old"; 7="Old"; 8="Young"; 9="Young"; 10="Young"; 11="Young"; 12="Young"; 13="Old"}

# Write the "Old" cells first so the shared-strings table ends up with
# "Old" before "Young" (rows 6, 7 and 13 are the first "Old" cages).
foreach ($r in 2..13) {
    if ($age[$r] -eq "Old") {
        $ws3.Cells.Item($r, 4).Value = $age[$r]
    }
}

foreach ($r in 2..13) {
    $ws3.Cells.Item($r, 3).Value = $grpSize[$r]
    if ($age[$r] -ne "Old") {
        $ws3.Cells.Item($r, 4).Value = $age[$r]
    }
}

# --- Rename "sexes" to "Sheet3" ---
$ws3.Name = "Sheet3"

# --- Keep the existing selections on Sheet1 / Sheet2 ---
$ws1.Activate()
$ws1.Range("H38").Select()

$ws2.Activate()
$ws2.Range("E79").Select()

# --- Sheet3 becomes the active tab with a new selection ---
$ws3.Activate()
$ws3.Range("F41").Select()
